# Rename worksheet "strategy_id-5008" -> "strategy_id-5007" and add a new
# worksheet "strategy_id-5009" (an exact duplicate of it, as the new
# calibrated strategy 5009 reuses the same template/data as 5007) placed
# right after it as the last sheet in the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("strategy_id-5008")

# Duplicate the sheet (copies all data, styles, sheetPr, page margins, etc.)
# and place the new copy immediately after the last sheet in the workbook.
$ws.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Rename the original sheet first ...
$ws.Name = "strategy_id-5007"

# ... then rename the freshly created copy (now the last sheet).
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "strategy_id-5009"

# Restore the originally active sheet/tab (copying a sheet makes the new
# copy active, which would otherwise shift the workbook's active tab).
$wb.Worksheets.Item(1).Activate()
